# Auto-generated Excel COM-interop PowerShell script
# Applies the crypto price/volume update described by the diff for Wed Sep 20 2023.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Cells.Item(2, 4).Value = "'27.101.63"
$ws.Cells.Item(2, 5).Value = "  -0.31%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "'1.624.07"
$ws.Cells.Item(3, 5).Value = "  -1.09%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.00%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'213.98"
$ws.Cells.Item(5, 5).Value = "  -1.49%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "'0.522"
$ws.Cells.Item(6, 5).Value = "  +1.42%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  +0.00%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  -1.49%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "'0.0627"
$ws.Cells.Item(9, 5).Value = "  -0.04%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "'20.28"
$ws.Cells.Item(10, 5).Value = "  +1.11%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  -0.06%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "'1.637.94"
$ws.Cells.Item(12, 5).Value = "  -0.46%  "

# Row 13
$ws.Cells.Item(13, 5).Value = "  -0.42%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "'0.543"
$ws.Cells.Item(14, 5).Value = "  -0.16%  "

# Row 15
$ws.Cells.Item(15, 2).Value = "WrappedBTC"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Cells.Item(15, 4).Value = "'27.101.72"
$ws.Cells.Item(15, 5).Value = "  -0.22%  "

# Row 16
$ws.Cells.Item(16, 2).Value = "Litecoin"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Cells.Item(16, 4).Value = "'64.61"
$ws.Cells.Item(16, 5).Value = "  -4.10%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "'0.0₃0742"
$ws.Cells.Item(17, 5).Value = "  +0.27%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "'216.70"
$ws.Cells.Item(18, 5).Value = "  -0.96%  "

# Row 19
$ws.Cells.Item(19, 5).Value = "  +0.02%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "'6.95"
$ws.Cells.Item(20, 5).Value = "  +1.16%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "'2.42"
$ws.Cells.Item(22, 5).Value = "  -6.45%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "'9.06"
$ws.Cells.Item(23, 5).Value = "  -1.71%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "'148.02"
$ws.Cells.Item(24, 5).Value = "  +0.19%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  +0.02%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "'7.29"
$ws.Cells.Item(26, 5).Value = "  -3.34%  "

# Row 27
$ws.Cells.Item(27, 5).Value = "  -0.78%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "'15.61"
$ws.Cells.Item(28, 5).Value = "  -1.12%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "'0.0507"
$ws.Cells.Item(29, 5).Value = "  -0.62%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "'3.36"
$ws.Cells.Item(31, 5).Value = "  -0.47%  "

# Row 32
$ws.Cells.Item(32, 5).Value = "  -1.19%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "'1.349.48"
$ws.Cells.Item(33, 5).Value = "  +6.05%  "

# Row 34
$ws.Cells.Item(34, 5).Value = "  -0.04%  "

# Row 35
$ws.Cells.Item(35, 5).Value = "  -0.44%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  -0.02%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "'0.553"
$ws.Cells.Item(37, 5).Value = "  +1.34%  "

# Row 38
$ws.Cells.Item(38, 4).Value = "'0.856"
$ws.Cells.Item(38, 5).Value = "  +0.10%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "  -0.06%  "

# Row 40
$ws.Cells.Item(40, 5).Value = "  -0.82%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  +0.18%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "'65.56"
$ws.Cells.Item(42, 5).Value = "  +6.13%  "

# Row 43
$ws.Cells.Item(43, 4).Value = "'5.23"
$ws.Cells.Item(43, 5).Value = "  -1.44%  "

# Row 44
$ws.Cells.Item(44, 4).Value = "'1.762.79"

# Row 45
$ws.Cells.Item(45, 4).Value = "'90.63"
$ws.Cells.Item(45, 5).Value = "  -1.35%  "

# Row 46
$ws.Cells.Item(46, 5).Value = "  +0.68%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "'0.859"
$ws.Cells.Item(47, 5).Value = "  +29.23%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "'0.0₆0105"
$ws.Cells.Item(48, 5).Value = "  -1.59%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "'0.0513"
$ws.Cells.Item(49, 5).Value = "  -0.11%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "'0.0992"
$ws.Cells.Item(50, 5).Value = "  +1.67%  "

# Row 51
$ws.Cells.Item(51, 4).Value = "'7.62"
$ws.Cells.Item(51, 5).Value = "  -0.94%  "
